# mher: Test for getdents.
#
# Row 5 on the sheet (ID 5, task "Improve file hiding") gets its
# description clarified and is now marked Done.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "Improve overall[, get rid of used syscall in kernel space]"
$ws.Range("G9").Value = "Yes"
$ws.Range("G9").Style = "Gut"

# Clear the "Done? = No" autofilter criterion so every row shows again.
$ws.ShowAllData()

# Restore the cursor to where the author left it.
$ws.Range("D13").Select() | Out-Null
